# Weekly update: a new price record for "Repollo" (Crespo record, Primera)
# from Terminal Hortofrutícola Agro Chillán is inserted as the new first
# data row of the table (row 271), pushing the existing rows 271-303 down
# to 272-304 (table grows from A1:R303 to A1:R304).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 271; everything below (old rows
# 271-303) shifts down to 272-304, carrying its data/formatting along.
$ws.Rows("271").Insert()

# Populate the newly inserted row 271 with the new record's values.
$ws.Cells.Item(271, 1).Value  = 7
$ws.Cells.Item(271, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(271, 3).Value  = "Ñuble"
$ws.Cells.Item(271, 4).Value  = 44918
$ws.Cells.Item(271, 5).Value  = 16
$ws.Cells.Item(271, 6).Value  = 100112006
$ws.Cells.Item(271, 7).Value  = "Repollo"
$ws.Cells.Item(271, 8).Value  = "Crespo record"
$ws.Cells.Item(271, 9).Value  = "Primera"
$ws.Cells.Item(271, 10).Value = 500
$ws.Cells.Item(271, 11).Value = 1300
$ws.Cells.Item(271, 12).Value = 1400
$ws.Cells.Item(271, 13).Value = 1350
$ws.Cells.Item(271, 14).Value = "`$/unidad"
$ws.Cells.Item(271, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(271, 16).Value = 1350
$ws.Cells.Item(271, 17).Value = 1
$ws.Cells.Item(271, 18).Value = "Hortaliza"
